$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.373.16"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "1.827.82"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3793"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07426"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8798"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").Value = "1.824.52"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.728"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.447"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07090"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008827"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Value = "27.366.95"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.342"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").Value = "2.045.67"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.276"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08957"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8001"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.198"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.557"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.932"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01978"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05256"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.336"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5349"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.374"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.64%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1709"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.691"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5133"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.689"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9992"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06390"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "
